$d = $word.ActiveDocument

# Locate the paragraph containing the copyright/footer notice that needs to
# be removed, along with the empty "page break before" paragraph that
# immediately precedes it.
$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $copyrightPara = $p
        break
    }
}

if ($copyrightPara -ne $null) {
    $prevPara = $copyrightPara.Previous()

    # Delete the copyright paragraph first, then the empty page-break
    # paragraph before it. Deleting in this order (last-to-first) keeps
    # each Range valid and preserves the paragraph mark/formatting of the
    # paragraphs that should remain untouched.
    $copyrightPara.Range.Delete()
    $prevPara.Range.Delete()
}
